# Auto-generated edit script applying the Spriggan_Profits diff
# Updates market-price / profit values across the ALC, ARM, BSM, CRP, CUL, GSM, LTW sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 6474.625
$ws.Range("I9").Value = 9277.909
$ws.Range("J9").Value = 307.4
$ws.Range("K9").Value = 9277.909
$ws.Range("L9").Value = 307.4
$ws.Range("M9").Value = -9108.909
$ws.Range("N9").Value = -645.4
$ws.Range("H51").Value = 8000
$ws.Range("I51").Value = 8000
$ws.Range("K51").Value = 8000
$ws.Range("M51").Value = -7516
$ws.Range("H55").Value = 49.5
$ws.Range("I55").Value = 50
$ws.Range("J55").Value = 49
$ws.Range("K55").Value = 50
$ws.Range("L55").Value = 49
$ws.Range("M55").Value = 164
$ws.Range("N55").Value = -477
$ws.Range("H61").Value = 1667081.6
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()
$ws.Range("H100").Value = 4266.222
$ws.Range("I100").Value = 4159.2
$ws.Range("K100").Value = 4159.2
$ws.Range("M100").Value = -3618.2
$ws.Range("H132").Value = 2541.5264
$ws.Range("I132").Value = 2541.5264
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 7624.5792
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -5094.5792
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H36").Value = 9206
$ws.Range("I36").Value = 3941.6667
$ws.Range("J36").Value = 24999
$ws.Range("K36").Value = 3941.6667
$ws.Range("L36").Value = 24999
$ws.Range("M36").Value = -3595.6667
$ws.Range("N36").Value = -25691
$ws.Range("H45").Value = 6232.8823
$ws.Range("I45").Value = 6632.3076
$ws.Range("J45").Value = 4934.75
$ws.Range("K45").Value = 6632.3076
$ws.Range("L45").Value = 4934.75
$ws.Range("M45").Value = -6255.3076
$ws.Range("N45").Value = -5688.75
$ws.Range("H74").Value = 29414704
$ws.Range("J74").Value = 5984.25
$ws.Range("L74").Value = 5984.25
$ws.Range("N74").Value = -7732.25
$ws.Range("H77").Value = 29414704
$ws.Range("J77").Value = 5984.25
$ws.Range("L77").Value = 29921.25
$ws.Range("N77").Value = -38657.25
$ws.Range("H122").Value = 3823.2307
$ws.Range("I122").Value = 3539.2354
$ws.Range("J122").Value = 4359.6665
$ws.Range("K122").Value = 10617.7062
$ws.Range("L122").Value = 13078.9995
$ws.Range("M122").Value = -8167.706200000001
$ws.Range("N122").Value = -17978.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3477.1667
$ws.Range("I86").Value = 3286.818
$ws.Range("J86").Value = 3776.2856
$ws.Range("K86").Value = 3286.818
$ws.Range("L86").Value = 3776.2856
$ws.Range("M86").Value = -2163.818
$ws.Range("N86").Value = -6022.2856
$ws.Range("H89").Value = 3477.1667
$ws.Range("I89").Value = 3286.818
$ws.Range("J89").Value = 3776.2856
$ws.Range("K89").Value = 16434.09
$ws.Range("L89").Value = 18881.428
$ws.Range("M89").Value = -10818.09
$ws.Range("N89").Value = -30113.428

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H29").Value = 7500
$ws.Range("J29").Value = 7500
$ws.Range("L29").Value = 7500
$ws.Range("N29").Value = -8086
$ws.Range("H62").Value = 2728.6667
$ws.Range("I62").Value = 1500
$ws.Range("J62").Value = 2974.4
$ws.Range("K62").Value = 1500
$ws.Range("L62").Value = 2974.4
$ws.Range("M62").Value = -876
$ws.Range("N62").Value = -4222.4
$ws.Range("H65").Value = 2728.6667
$ws.Range("I65").Value = 1500
$ws.Range("J65").Value = 2974.4
$ws.Range("K65").Value = 7500
$ws.Range("L65").Value = 14872
$ws.Range("M65").Value = -4380
$ws.Range("N65").Value = -21112
$ws.Range("H68").Value = 79259.89
$ws.Range("I68").Value = 39985.668
$ws.Range("J68").Value = 98897
$ws.Range("K68").Value = 39985.668
$ws.Range("L68").Value = 98897
$ws.Range("M68").Value = -39236.668
$ws.Range("N68").Value = -100395
$ws.Range("H71").Value = 79259.89
$ws.Range("I71").Value = 39985.668
$ws.Range("J71").Value = 98897
$ws.Range("K71").Value = 119957.004
$ws.Range("L71").Value = 296691
$ws.Range("M71").Value = -116213.004
$ws.Range("N71").Value = -304179
$ws.Range("H94").Value = 2276.7
$ws.Range("J94").Value = 1832.6666
$ws.Range("L94").Value = 1832.6666
$ws.Range("N94").Value = -2734.6666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 320.9
$ws.Range("I14").Value = 320.9
$ws.Range("K14").Value = 962.6999999999999
$ws.Range("M14").Value = -789.6999999999999
$ws.Range("H107").Value = 728.16
$ws.Range("I107").Value = 309.1111
$ws.Range("J107").Value = 963.875
$ws.Range("K107").Value = 927.3333
$ws.Range("L107").Value = 2891.625
$ws.Range("M107").Value = 992.6667
$ws.Range("N107").Value = -6731.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 203579.84
$ws.Range("J21").Value = 28999.555
$ws.Range("L21").Value = 28999.555
$ws.Range("N21").Value = -29345.555
$ws.Range("H30").Value = 203579.84
$ws.Range("J30").Value = 28999.555
$ws.Range("L30").Value = 28999.555
$ws.Range("N30").Value = -29209.555
$ws.Range("H97").Value = 1650.1578
$ws.Range("I97").Value = 1608.2307
$ws.Range("K97").Value = 1608.2307
$ws.Range("M97").Value = -1112.2307
$ws.Range("H104").Value = 25192.5
$ws.Range("J104").Value = 25192.5
$ws.Range("L104").Value = 25192.5
$ws.Range("N104").Value = -32180.5
$ws.Range("H107").Value = 8000.6665
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 8000.6665
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 8000.6665
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = -11840.6665
$ws.Range("H113").Value = 63022.59
$ws.Range("I113").Value = 81629.62
$ws.Range("K113").Value = 81629.62
$ws.Range("M113").Value = -79459.62

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 1000000
$ws.Range("J43").Value = 1000000
$ws.Range("L43").Value = 1000000
$ws.Range("N43").Value = -1000386
$ws.Range("H122").Value = 6017.227
$ws.Range("I122").Value = 4419
$ws.Range("J122").Value = 21999.5
$ws.Range("K122").Value = 13257
$ws.Range("L122").Value = 65998.5
$ws.Range("M122").Value = -10807
$ws.Range("N122").Value = -70898.5
